# Generate Report for handoff
# Updates the "Latest Handoff Datetime" for the 71d4124f-... row on both the
# zh-cn and de-de localization-status sheets to reflect a fresh handoff run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-15 07:16:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-15 07:16:22"
